$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("UF1")
$ws.Range("F1").Value = "ejercicio1"
